$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matching the inline/shared-string cell type
# used throughout this sheet), then drop back to the default "Normal" style so
# no stray cell-level style reference is left behind.
$textCells = @("D5", "D6", "D9", "D10", "D12", "D17", "D18", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D46", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "71.476.17"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "3.977.11"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "538.09"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("D6").Value = "150.85"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("D7").Value = "3.970.97"
$ws.Range("E7").Value = "  -2.19%  "
$ws.Range("E8").Value = "  -4.12%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "0.748"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -7.50%  "
$ws.Range("D12").Value = "55.49"
$ws.Range("E12").Value = "  +13.54%  "
$ws.Range("E13").Value = "  -5.05%  "
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").Value = "4.609.19"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "3.966.18"
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("D17").Value = "14.06"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").Value = "20.66"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("E20").Value = "  -4.88%  "
$ws.Range("D21").Value = "71.336.94"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").Value = "434.31"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "97.51"
$ws.Range("E24").Value = "  -6.09%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "14.68"
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "4.24"
$ws.Range("E26").Value = "  +4.19%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "11.39"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "4.09"
$ws.Range("E28").Value = "  +23.75%  "
$ws.Range("D29").Value = "10.89"
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("D30").Value = "5.91"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").Value = "36.86"
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("D32").Value = "7.99"
$ws.Range("E32").Value = "  +19.82%  "
$ws.Range("D33").Value = "51.92"
$ws.Range("E33").Value = "  +21.96%  "
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D36").Value = "678.41"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "65.87"
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("D38").Value = "0.447"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").Value = "0.0₃0828"
$ws.Range("E39").Value = "  -8.54%  "
$ws.Range("E40").Value = "  -3.28%  "
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").Value = "10.46"
$ws.Range("E46").Value = "  +5.07%  "
$ws.Range("E47").Value = "  -4.89%  "
$ws.Range("D48").Value = "2.68"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").Value = "3.35"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("E51").Value = "  -9.23%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
